# Leave Card update (5/18/2023 4:35 PM):
# Two new SL(3-0-0) absence rows are inserted into the leave table (Sheet1 / Table1)
# right after the existing 5/2023 monthly row (row 650), pushing every subsequent
# monthly row down by two and growing the table by two rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# --- 1. Insert two new blank rows at row 651 (shifts 651..762 down to 653..764) ---
$ws.Rows.Item(651).Insert()
$ws.Rows.Item(651).Insert()

# --- 2. Grow the table definition to match the new bottom of the data (K764) ---
$tbl.Resize($ws.Range("A8:K764"))

# --- 3. Re-assert the calculated "EARNED " column formula for the two rows that
#        rolled off the end of the originally-sized table so they keep the full
#        structured reference (rather than an out-of-table literal reference). ---
$ws.Range("G763").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G764").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 4. Copy the formatting of the row above (650, the last real monthly row)
#        down onto the two freshly-inserted rows - this is what Excel does natively
#        when a row is inserted in the middle of a table. ---
$ws.Range("A650:K650").Copy()
$ws.Range("A651:K652").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Populate the two new SL(3-0-0) leave rows ---
# Row 651: SL(3-0-0), 3 days, covering 5/11,12,15/2023
$ws.Range("B651").Value = "SL(3-0-0)"
$ws.Range("H651").Value = 3
$ws.Range("K651").Value = "5/11,12,15/2023"

# Row 652: SL(3-0-0), 3 days, covering 5/8-10/2023
$ws.Range("B652").Value = "SL(3-0-0)"
$ws.Range("H652").Value = 3
$ws.Range("K652").Value = "5/8-10/2023"

# Re-apply the calculated "EARNED " column formula on the two new rows too, since
# the PasteSpecial(formats only) step above cleared the formula that was copied in
# at Insert() time.
$ws.Range("G651").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G652").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 6. Leave the selection on the new blank particulars cell, matching where the
#        editor's cursor ended up after the insert. ---
$ws.Activate()
$ws.Range("B653").Select()

Write-Output "Leave Card update applied."
